$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (rename "Publication Identifier" / "filter" labels to the new Study-level labels) ---
$ws.Range("A1").Value = 'Study Identifier'
$ws.Range("H1").Value = 'Population filter 1'
$ws.Range("J1").Value = 'Population filter 2'

# --- Write the new data rows 2-11, column by column (A..L) so the cell-write order matches the
#     original authoring order (this also keeps the shared-string table build-up deterministic) ---
# Column A
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 6
$ws.Range("A7").Value = 12
$ws.Range("A8").Value = 13
$ws.Range("A9").Value = 14
$ws.Range("A10").Value = 16
$ws.Range("A11").Value = 19

# Column B
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 6
$ws.Range("B7").Value = '12, 13'
$ws.Range("B8").Value = 14
$ws.Range("B9").Value = 15
$ws.Range("B10").Value = 17
$ws.Range("B11").Value = 20

# Column C
$ws.Range("C2").Value = 'Original'
$ws.Range("C3").Value = 'Original'
$ws.Range("C4").Value = 'Original'
$ws.Range("C5").Value = 'Original'
$ws.Range("C6").Value = 'Original'
$ws.Range("C7").Value = 'Original'
$ws.Range("C8").Value = 'Original'
$ws.Range("C9").Value = 'Original'
$ws.Range("C10").Value = 'Original'
$ws.Range("C11").Value = 'Original'

# Column D
$ws.Range("D2").Value = 'Bertranou_2017'
$ws.Range("D3").Value = 'Horgan_2010'
$ws.Range("D4").Value = 'Wu_2018'
$ws.Range("D5").Value = 'Rui_2020'
$ws.Range("D6").Value = 'Guan_2019'
$ws.Range("D7").Value = 'TA416, TA653: Osimertinib for treating EGFR T790M mutation-positive advanced non-small-cell lung cancer'
$ws.Range("D8").Value = 'NICE_TA584'
$ws.Range("D9").Value = 'SMC_1214/17'
$ws.Range("D10").Value = 'pCODR 2020'
$ws.Range("D11").Value = 'PBAC_2019'

# Column E
$ws.Range("E2").Value = 'Cost-effectiveness of osimertinib in the UK for advanced EGFR-T790M non-small cell lung cancer'
$ws.Range("E3").Value = 'An economic analysis of the INTEREST trial, a randomized trial of docetaxel versus gefitinib as second-/third-line therapy in advanced non-small-cell lung cancer'
$ws.Range("E4").Value = 'Cost-Effectiveness of Osimertinib for EGFR Mutation–Positive Non–Small Cell Lung Cancer after Progression following First-Line EGFR TKI Therapy'
$ws.Range("E5").Value = 'Cost-effectiveness of Osimertinib vs Docetaxel-bevacizumab in Third-line Treatment in EGFR T790M Resistance Mutation Advanced Non-Small Cell Lung Cancer in China'
$ws.Range("E6").Value = 'Cost-effectiveness of Osimertinib as a Second-line Treatment in Patients With EGFR-mutated Advanced Non-Small Cell Lung Cancer in China'
$ws.Range("E7").Value = 'TA416, TA653: Osimertinib for treating EGFR T790M mutation-positive advanced non-small-cell lung cancer'
$ws.Range("E8").Value = 'TA584: Atezolizumab in combination for treating metastatic non-squamous non-small-cell lung cancer'
$ws.Range("E9").Value = 'Osimertinib 40mg and 80mg film-coated tablets (Tagrisso®)'
$ws.Range("E10").Value = 'Tecentriq & Avastin Non-Squamous Non-Small Cell Lung Cancer'
$ws.Range("E11").Value = 'Atezolizumab and Bevacizumab: atezolizumab: Solution concentrate for I.V. infusion 1200 mg in 20 mL; bevacizumab: Solution for I.V. infusion 100 mg in 4 mL, Solution for I.V. infusion 400 mg in 16 mL; Tecentriq® and Avastin®'

# Column F
$ws.Range("F2").Value = 'NR'
$ws.Range("F3").Value = 'NR'
$ws.Range("F4").Value = 'NR'
$ws.Range("F5").Value = 'NR'
$ws.Range("F6").Value = 'NR'
$ws.Range("F7").Value = '1) Tseng et al. 2014
2) NEJ002 (Miyauchi et al. 2015)
3) Park et al. 2015
4) Halmos et al. 2015 (RCT)
5) Zhou et al. 2014
6) Kashara et al. 2015 (RCT)
7) IMPRESS (Soria et al. 2015; RCT)
8) Wu et al. 2010
9) Kim et al. 2013
10) AURA 2/AURA extension (single-arm)'
$ws.Range("F8").Value = 'IMpower150 (NCT02366143)
Ongoing at time of submission'
$ws.Range("F9").Value = 'NR'
$ws.Range("F10").Value = 'IMpower150 (NCT02366143)'
$ws.Range("F11").Value = 'NR'

# Column G
$ws.Range("G2").Value = 'Patients with EGFR-T790M mutation positive NSCLC who have progressed on or after EGFR-TKI therapy.'
$ws.Range("G3").Value = 'Patients with locally advanced or metastatic NSCLC that had progressed/recurred following one to two prior chemotherapy regimens including platinum'
$ws.Range("G4").Value = 'Patients with EGFR T790M-positive advanced NSCLC after the failure of first-line therapy with first-generation EGFR TKIs'
$ws.Range("G5").Value = 'Advanced NSCLC with acquired EGFR T790M resistance mutation previously treated with gefitinib or erlotinib in first-line and platinum-based chemotherapy in second-line'
$ws.Range("G6").Value = 'Advanced, EGFR mutation-positive NSCLC after failure of EGFR TKI'
$ws.Range("G7").Value = 'Treating epidermal growth factor receptor (EGFR) T790M mutation-positive locally advanced or metastatic non-small-cell lung cancer (NSCLC)'
$ws.Range("G8").Value = 'Metastatic non-squamous non-small-cell lung cancer'
$ws.Range("G9").Value = 'Treatment of adult patients with locally advanced or metastatic epidermal growth factor receptor (EGFR) T790M mutation-positive non-small-cell lung cancer (NSCLC).'
$ws.Range("G10").Value = 'Metastatic EGFR and/or ALK-positive NSCLC in patients who have progressed on treatment with targeted therapies.'
$ws.Range("G11").Value = 'Locally advanced or metastatic EGFR/ALK NSCLC patients, who have disease progression on or after treatment with an EGFR/ALK TKI'

# Column H
$ws.Range("H2").Value = 'EGFRm ITT'
$ws.Range("H3").Value = 'EGFRm subgroup'
$ws.Range("H4").Value = 'EGFRm subgroup'
$ws.Range("H5").Value = 'EGFRm ITT'
$ws.Range("H6").Value = 'EGFRm ITT'
$ws.Range("H7").Value = 'EGFRm ITT'
$ws.Range("H8").Value = 'EGFRm (mixed)'
$ws.Range("H9").Value = 'EGFRm ITT'
$ws.Range("H10").Value = 'EGFRm (mixed)'
$ws.Range("H11").Value = 'EGFRm (mixed)'

# Column I
$ws.Range("I2").Value = 'NR'
$ws.Range("I3").Value = 'Patients with locally advanced or metastatic NSCLC that had progressed/recurred following one to two prior chemotherapy regimens including platinum'
$ws.Range("I4").Value = 'NR'
$ws.Range("I5").Value = 'NR'
$ws.Range("I6").Value = 'Patients with advanced NSCLC, a mean age of 65 years, and a World Health Organization (WHO) performance status of 0-2.'
$ws.Range("I7").Value = 'NR'
$ws.Range("I8").Value = 'NR'
$ws.Range("I9").Value = 'NR'
$ws.Range("I10").Value = 'NR'
$ws.Range("I11").Value = 'NR'

# Column J
$ws.Range("J2").Value = '2+ Line'
$ws.Range("J3").Value = '2+ Line'
$ws.Range("J4").Value = '2+ Line'
$ws.Range("J5").Value = '2+ Line'
$ws.Range("J6").Value = '2+ Line'
$ws.Range("J7").Value = '2+ Line'
$ws.Range("J8").Value = '2+ Line'
$ws.Range("J9").Value = '2+ Line'
$ws.Range("J10").Value = '2+ Line'
$ws.Range("J11").Value = '2+ Line'

# Column K
$ws.Range("K2").Value = '2015 / UK'
$ws.Range("K3").Value = '2008 / Canada'
$ws.Range("K4").Value = '2017 / China, USA'
$ws.Range("K5").Value = '2019 / China'
$ws.Range("K6").Value = '2018 / China'
$ws.Range("K7").Value = '2016 / UK'
$ws.Range("K8").Value = '2019 / UK'
$ws.Range("K9").Value = '2017 / UK'
$ws.Range("K10").Value = '2020 / Canada'
$ws.Range("K11").Value = '2019 / Australia'

# Column L
$ws.Range("L2").Value = 'Progression-free
Progressed disease
Death'
$ws.Range("L3").Value = 'Stable disease
Responsive disease
Progressive disease'
$ws.Range("L4").Value = 'Progression-free
Progressed-survival
Death'
$ws.Range("L5").Value = 'Progression-free survival
Postprogression survival
Death'
$ws.Range("L6").Value = 'Progression-free
Progression
Death'
$ws.Range("L7").Value = 'TA416:
Base case
Progression-free: 0.815
Post-progression: 0.678
Second-line only population (subgroup)
Progression-free: 0.853
Post-progression: 0.726
≥Third-line population (subgroup)
Progression-free: 0.798
Post-progression: 0.659
Updated base case (osimertinib)
Progression-free: 0.806 (adjusted dataset)
Progression-free: 0.805 (unadjusted dataset)
Progressed disease: 0.715
Updated base case (PDC)
Progression-free: 0.779 (adjusted dataset)
Progression-free: 0.783 (unadjusted dataset)
Progressed disease: 0.715
TA653:
Model A
PF: 0.831
Stable disease: 0.751
PD: 0.715
Model B
PF: 0.836
Stable disease: 0.797
PD: 0.717
Post-technical engagement, the company updated their base case and used the following utility values:
Osimertinib
Response: 0.831
Stable disease: 0.751
Progressed disease: 0.715
PDC
Response: 0.67
Stable disease: 0.67
Progressed disease: 0.64'
$ws.Range("L8").Value = '≤ 5 weeks before death: 0.52
> 5 and ≤ 15 weeks before death:  0.59
> 15 and ≤ 30 weeks before death: 0.70 
> 30 weeks before death: 0.73'
$ws.Range("L9").Value = 'Base case:
Progression free state: 0.831
Stable disease: 0.751
Progressed disease: 0.715'
$ws.Range("L10").Value = 'Health state utility values in the model were based on a patient’s proximity to death, with changes in patient utility occurring independently of progression, as of 30 weeks from death.'
$ws.Range("L11").Value = 'Increment LYG: 0.785                                                         Increment QALYS: 0.565'

# --- Some of the new cells contain embedded line breaks; let Excel recompute the row heights
#     back to automatic (the rows keep their default height in the saved workbook) ---
$ws.Range("A2:A11").EntireRow.AutoFit()

# --- Restore the author's final selection ---
$ws.Range("C5").Select() | Out-Null
